$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header) ---
# A1 already contains "Year" (shared string 0) with the bold/bordered style.
$ws.Range("B1").Value = "Flexible working hours"
$ws.Range("C1").Value = "Full-time"
$ws.Range("D1").Value = "Of which: female"
$ws.Range("E1").Value = "Of which: male"
$ws.Range("F1").Value = "Part-time"
$ws.Range("G1").Value = "Of which: female"
$ws.Range("H1").Value = "Of which: male"
$ws.Range("I1").Value = "Virtual offices"
$ws.Range("J1").Value = "Sabbatical"
$ws.Range("K1").Value = "Semi-retirement (Altersteilzeit)"

# --- Row 2 (new data row, year 2019) ---
$ws.Range("A2").Value = 2019
$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 78.90000000000001
$ws.Range("D2").Value = 58.4
$ws.Range("E2").Value = 93.90000000000001
$ws.Range("F2").Value = 21.1
$ws.Range("G2").Value = 41.6
$ws.Range("H2").Value = 6.1
$ws.Range("I2").Value = 13.2
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1.76

# --- Row 3 (previously blank placeholder row, year 2018) ---
$ws.Range("A3").Value = 2018
$ws.Range("B3").Value = 100
$ws.Range("C3").Value = 79.90000000000001
$ws.Range("D3").Value = 59.9
$ws.Range("E3").Value = 94.59999999999999
$ws.Range("F3").Value = 20.1
$ws.Range("G3").Value = 40.1
$ws.Range("H3").Value = 5.4
$ws.Range("I3").Value = 12.8
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1.68

# --- Row 4 (brand new row, year 2017) ---
# Copy the formatting of A3 (bold, centered, bordered) onto A4 before
# assigning its value, so the new row's first cell keeps the same style.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A4").Value = 2017
$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 80.8
$ws.Range("F4").Value = 19.2
$ws.Range("I4").Value = 12.9
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1.8
